# Implement mid-period revenue and expense recognition in cost income rule
# Adds two new rows (Project E / Project F) to the "costs" sheet, reusing
# the existing date-formatted style from the row above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("costs")

# Copy the formatting of the last existing data row (row 7) down onto the
# two new rows so the date columns (A, D, E) pick up the same number
# format / style index instead of minting a brand-new style.
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122)

# Row 8 - Project E (revenue recognised mid-period)
$ws.Range("A8").Value = 46174
$ws.Range("B8").Value = "Project E"
$ws.Range("C8").Value = 5000
$ws.Range("D8").Value = 45717
$ws.Range("E8").Value = 45961

# Row 9 - Project F (expense recognised mid-period)
$ws.Range("A9").Value = 46174
$ws.Range("B9").Value = "Project F"
$ws.Range("C9").Value = -4000
$ws.Range("D9").Value = 45717
$ws.Range("E9").Value = 45961

# Match the author's final selection (moved down one row to A8).
[void]$ws.Range("A8").Select()
